$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow the data refresh, then restore protection.
$ws.Unprotect("D382")

# Update the "as of" date in the confidentiality banner (shared string).
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."
$ws.Rows(41).AutoFit()

# Refresh the Weight (D) and Percent Change (E) columns with the latest snapshot values.
$ws.Range("D2").Value = 0.02905729654914223
$ws.Range("E2").Value = -0.01697825257535301
$ws.Range("D3").Value = 0.02943423210147754
$ws.Range("E3").Value = -0.01888309430682311
$ws.Range("D4").Value = 0.03007894539913702
$ws.Range("E4").Value = -0.02424087777494244
$ws.Range("D5").Value = 0.06689113662058943
$ws.Range("E5").Value = -0.0160733067729083
$ws.Range("D6").Value = 0.0144207696606109
$ws.Range("E6").Value = 0.002956830277942135
$ws.Range("D7").Value = 0.01412655525890121
$ws.Range("E7").Value = -0.04029580440688185
$ws.Range("D8").Value = 0.03134534651953957
$ws.Range("E8").Value = 0.001754815540319976
$ws.Range("D9").Value = 0.03120804646540839
$ws.Range("E9").Value = -0.01420959147424505
$ws.Range("D10").Value = 0.03422417048061629
$ws.Range("E10").Value = -0.02668709936646185
$ws.Range("D11").Value = 0.03040258124101767
$ws.Range("E11").Value = -0.03399671813859551
$ws.Range("D12").Value = 0.01488298474967365
$ws.Range("E12").Value = -0.01432500573000228
$ws.Range("D13").Value = 0.01609992373877428
$ws.Range("E13").Value = -0.02651093807934735
$ws.Range("D14").Value = 0.01474142071870609
$ws.Range("E14").Value = 0.00138840680319352
$ws.Range("D15").Value = 0.007844011788190363
$ws.Range("E15").Value = 0.01032833224614049
$ws.Range("D16").Value = 0.006617692050049708
$ws.Range("E16").Value = 0.01778350515463933
$ws.Range("D17").Value = 0.03208642569370106
$ws.Range("E17").Value = -0.02544186046511621
$ws.Range("D18").Value = 0.03244801092942543
$ws.Range("E18").Value = -0.02406764960971386
$ws.Range("D19").Value = 0.03188068881134611
$ws.Range("E19").Value = -0.01304042531848726
$ws.Range("D20").Value = 0.03098098969887157
$ws.Range("E20").Value = -0.02921240064687081
$ws.Range("D21").Value = 0.04376886942999372
$ws.Range("E21").Value = -0.003848102252357655
$ws.Range("D22").Value = 0.02829319190006432
$ws.Range("E22").Value = -0.00940410525363955
$ws.Range("D23").Value = 0.03091170007528052
$ws.Range("E23").Value = 0.009552382922960234
$ws.Range("D24").Value = 0.02992352344345125
$ws.Range("E24").Value = 0.02295607566527735
$ws.Range("D25").Value = 0.01430990626286522
$ws.Range("E25").Value = -0.02765196662693692
$ws.Range("D26").Value = 0.01334091752679963
$ws.Range("E26").Value = 0.0002397123451858274
$ws.Range("D27").Value = 0.0305953129940217
$ws.Range("E27").Value = 0.001226429556952402
$ws.Range("D28").Value = 0.03213844621110481
$ws.Range("E28").Value = -0.003343416653399101
$ws.Range("D29").Value = 0.03039106850355947
$ws.Range("E29").Value = -0.008923310042932897
$ws.Range("D30").Value = 0.02923382519016804
$ws.Range("E30").Value = -0.02887981330221712
$ws.Range("D31").Value = 0.03344002514040744
$ws.Range("E31").Value = -0.03272574259319483
$ws.Range("D32").Value = 0.03270598152802596
$ws.Range("E32").Value = -0.00466083464792777
$ws.Range("D33").Value = 0.02939329792384836
$ws.Range("E33").Value = -0.05161458786665507
$ws.Range("D34").Value = 0.03142550928406337
$ws.Range("E34").Value = 0.005698778833107321
$ws.Range("D35").Value = 0.03126752894227577
$ws.Range("E35").Value = 0.000231830300220448
$ws.Range("D36").Value = 0.02867822900838874
$ws.Range("E36").Value = 0.007761273919443257
$ws.Range("D37").Value = 0.03141143816050334
$ws.Range("E37").Value = -0.006196804539345901
$ws.Range("E38").Value = -0.01263416336617074

# Restore the original sheet protection.
$ws.Protect("D382")
